$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $Text) {
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "25.767.21"
Set-TextValue $ws "E2" "  +0.09%  "
Set-TextValue $ws "D3" "1.628.59"
Set-TextValue $ws "E3" "  -0.36%  "
Set-TextValue $ws "D5" "213.73"
Set-TextValue $ws "E5" "  -0.80%  "
Set-TextValue $ws "E6" "  -0.38%  "
Set-TextValue $ws "D7" "0.997"
Set-TextValue $ws "E7" "  -0.59%  "
Set-TextValue $ws "E8" "  -1.07%  "
Set-TextValue $ws "E9" "  -1.03%  "
Set-TextValue $ws "E10" "  +0.05%  "
Set-TextValue $ws "E11" "  +0.23%  "
Set-TextValue $ws "E12" "  +0.20%  "
Set-TextValue $ws "D13" "1.853.63"
Set-TextValue $ws "E13" "  -0.35%  "
Set-TextValue $ws "D14" "1.599.16"
Set-TextValue $ws "E14" "  -2.23%  "
Set-TextValue $ws "E15" "  -1.09%  "
Set-TextValue $ws "E16" "  -1.41%  "
Set-TextValue $ws "D17" "62.63"
Set-TextValue $ws "D18" "25.782.73"
Set-TextValue $ws "E18" "  +0.08%  "
Set-TextValue $ws "E19" "  -0.56%  "
Set-TextValue $ws "E20" "  -0.34%  "
Set-TextValue $ws "D21" "190.70"
Set-TextValue $ws "E21" "  -1.49%  "
Set-TextValue $ws "E22" "  -0.31%  "
Set-TextValue $ws "D23" "6.27"
Set-TextValue $ws "E23" "  -0.11%  "
Set-TextValue $ws "D24" "0.997"
Set-TextValue $ws "E24" "  -0.61%  "
Set-TextValue $ws "E25" "  -2.10%  "
Set-TextValue $ws "D26" "142.11"
Set-TextValue $ws "E26" "  +1.41%  "
Set-TextValue $ws "E27" "  +0.55%  "
Set-TextValue $ws "E28" "  -0.89%  "
Set-TextValue $ws "D29" "15.51"
Set-TextValue $ws "E29" "  +0.00%  "
Set-TextValue $ws "E30" "  -0.90%  "
Set-TextValue $ws "D31" "0.0493"
Set-TextValue $ws "E31" "  -0.20%  "
Set-TextValue $ws "E32" "  -0.66%  "
Set-TextValue $ws "E33" "  -0.90%  "
Set-TextValue $ws "E34" "  -0.42%  "
Set-TextValue $ws "E35" "  +0.06%  "
Set-TextValue $ws "E36" "  +0.26%  "
Set-TextValue $ws "D37" "1.140.21"
Set-TextValue $ws "E37" "  +1.53%  "
Set-TextValue $ws "E38" "  -1.00%  "
Set-TextValue $ws "E39" "  -1.36%  "
Set-TextValue $ws "E40" "  -0.26%  "
Set-TextValue $ws "D41" "0.995"
Set-TextValue $ws "E41" "  -0.73%  "
Set-TextValue $ws "D42" "5.59"
Set-TextValue $ws "E42" "  +0.54%  "
Set-TextValue $ws "D43" "100.36"
Set-TextValue $ws "E43" "  +0.73%  "
Set-TextValue $ws "D44" "0.798"
Set-TextValue $ws "E44" "  -0.38%  "
Set-TextValue $ws "D45" "1.764.62"
Set-TextValue $ws "E45" "  -0.26%  "
Set-TextValue $ws "B46" "BabyDogeCoin"
Set-TextValue $ws "C46" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D46" "0.0$([char]0x2086)0110"
Set-TextValue $ws "E46" "  -0.04%  "
Set-TextValue $ws "B47" "Aave"
Set-TextValue $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D47" "55.30"
Set-TextValue $ws "E47" "  +0.40%  "
Set-TextValue $ws "B48" "Cronos"
Set-TextValue $ws "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D48" "0.0511"
Set-TextValue $ws "E48" "  +2.11%  "
Set-TextValue $ws "B49" "RenderToken"
Set-TextValue $ws "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D49" "1.45"
Set-TextValue $ws "E49" "  +5.13%  "
Set-TextValue $ws "B50" "Mantle"
Set-TextValue $ws "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D50" "0.415"
Set-TextValue $ws "E50" "  -0.37%  "
Set-TextValue $ws "B51" "EnergySwap"
Set-TextValue $ws "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D51" "7.51"
Set-TextValue $ws "E51" "  -0.88%  "
